$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("Data")

# --- Update the Data (lookup) sheet ---------------------------------------
# Order matters here only insofar as it controls the order brand-new shared
# strings get appended to the shared-string table; it mirrors the order the
# original author appears to have typed the new values in.

# Age bracket column (C): replace "75+" with "75 and over" and add a new
# "No Response" option below it.
$data.Range("C16").Value = "75 and over"
$data.Range("C17").Value = "No Response"

# Gender column (B): add "Not listed" and "No Response" options.
$data.Range("B3").Value = "Not listed"
$data.Range("B4").Value = "No Response"

# Residence column (D): drop "Another Province/Country", reorder
# Calgary/Rocky View, and add a "No Response" option.
$data.Range("D5").Value = "No Response"

# Language column (E): add French and Chinese, drop the old
# "Another Province/Country"-era ordering, and shuffle remaining languages.
$data.Range("E2").Value = "French"
$data.Range("E3").Value = "Chinese"
$data.Range("E4").Value = "German"

# Income column (F): replace "I don't want to say" with "Prefer not to say"
# and add a "No Response" option.
$data.Range("F7").Value = "Prefer not to say"
$data.Range("F8").Value = "No Response"

# Remaining language / residence reshuffles (values already exist in the
# shared-string table by this point, so these only move which cell points at
# which string).
$data.Range("E5").Value = "Spanish"
$data.Range("E6").Value = "Punjabi"
$data.Range("E7").Value = "Tagalog (Pilipino)"
$data.Range("E8").Value = "Vietnamese"
$data.Range("E9").Value = "Other"
$data.Range("E10").Value = "No Response"

$data.Range("D2").Value = "Rocky View & County Area"
$data.Range("D3").Value = "Calgary"
$data.Range("D4").Value = "Other"

# --- Column widths on the Data sheet (widened to fit the new options) -----
$data.Columns.Item(2).AutoFit()
$data.Columns.Item(3).AutoFit()

# --- Switch the active tab from Input to Data ------------------------------
$data.Range("F9").Select()
$data.Activate()
